$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (becomes the data that used to be in row 5)
$ws.Range("A2").Value = 111525223
$ws.Range("B2").Value = 78107
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = "Vedskivlav"
$ws.Range("G2").Value = "Hertelidea botryosa"
$ws.Range("H2").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q2").Value = 404637.0659126193
$ws.Range("R2").Value = 6706784.214121711

# Row 4 (becomes the data that used to be in row 6)
$ws.Range("A4").Value = 111525235
$ws.Range("Q4").Value = 404485.2245768273
$ws.Range("R4").Value = 6706757.647421388

# Row 5 (becomes the data that used to be in row 7)
$ws.Range("A5").Value = 111525233
$ws.Range("B5").Value = 77515
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 404540.9329893424
$ws.Range("R5").Value = 6706716.233959051

# Row 6 (becomes the data that used to be in row 2)
$ws.Range("A6").Value = 111525238
$ws.Range("Q6").Value = 404495.4563026094
$ws.Range("R6").Value = 6706677.491168984

# Row 7 (becomes the data that used to be in row 4)
$ws.Range("A7").Value = 111525226
$ws.Range("Q7").Value = 404616.9589749529
$ws.Range("R7").Value = 6706770.937089294
